# End-of-week T/L, journal update
# Almost done with RecCalc.
#
# Adds Thursday/Friday clock-in/out times, corrects the Saturday
# clock-in/out times, annotates Saturday (L18) with a note about the
# Friday entries carrying onto Sunday, and nudges the selection the way
# the source workbook shows after the edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Thursday (row 16): add In / Out times -> 11:00 AM - 2:00 PM
$ws.Cells.Item(16, 3).Value = 0.458333333333333
$ws.Cells.Item(16, 4).Value = 0.583333333333333

# Friday (row 17): add In / Out times -> 12:00 PM - 3:00 PM
$ws.Cells.Item(17, 3).Value = 0.5
$ws.Cells.Item(17, 4).Value = 0.625

# Saturday (row 18): correct In / Out times -> 4:50 PM - 6:15 PM
$ws.Cells.Item(18, 3).Value = 0.701388888888889
$ws.Cells.Item(18, 4).Value = 0.760416666666667

# Saturday row annotation (new column L) explaining the Friday entries
# that carried onto Sunday because of a seminar clock-out.
$ws.Cells.Item(18, 12).Value = "// Multiple time entries on Friday carried onto Sunday to account for clocking out for seminar"

# Append the next two auto-numbered print-area names, matching the
# pattern already present in the workbook's defined names.
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# Leave the selection where the author's session ended up.
$ws.Range("V31").Select() | Out-Null
